$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEB-22")
$ws.Activate()

# --- Row 16: new entry (Task #9, 9-Feb-2022) ---
$ws.Range("A16").Value = 9
$ws.Range("B16").Value = 44601
$ws.Range("C16").Value = "RPA GSS"
$ws.Range("D16").Value = "1. Customization due to Multi-Factor Authentication  has been implemented at GRS-Details, tested and it is running smoothly"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = "Completed"

# B16 was not previously a date-formatted cell (unlike B17 and the rows above
# it) - pull the date number format across from B17 so it renders the same
# way Excel left it, reusing the existing style record instead of minting a
# new one.
$ws.Range("B17").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null

# --- Row 17: new entry (Task #10, 10-Feb-2022) ---
$ws.Range("A17").Value = 10
$ws.Range("B17").Value = 44602
$ws.Range("C17").Value = "RPA GSS"
$ws.Range("D17").Value = "1. A 12 daily tasks has been customised due to MFA issue, completed, tested and they are running smoothly"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = "Completed"

$excel.CutCopyMode = $false

# --- Update the saved view: scrolled down a bit, active cell on F17 ---
$win = $excel.ActiveWindow
$ws.Range("F17").Select() | Out-Null
$win.ScrollRow = 13
$win.ScrollColumn = 1 | Out-Null
